$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# Scroll the sheet back so column A/C is visible again (it had been
# scrolled over to show column Q), then clear the "Units" (M) column
# values for the data rows, leaving that range selected.
$ws.Range("C1").Select()
$ws.Range("M2:M6").Select()
$ws.Range("M2:M6").ClearContents()
